# PNAD 2009 - "tentativa" sheet correction
#
# The previous layout had a spurious header row ("grandes regiões e
# unidades da federação") at row 6 with no data of its own, followed by
# the real region rows (norte, rondônia, acre, ...) in rows 7-38.
#
# The fix removes that stray header row entirely: deleting the row
# shifts every region row (and its B:G data) up by one, the row 38
# previously holding "distrito federal" collapses into the new last
# row 37, the sheet dimension shrinks from A1:G38 to A1:G37, and the
# now-unused shared string for the removed header text drops out of
# sharedStrings.xml (count/uniqueCount 48 -> 47).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("6").Delete()
